# Rotate columns A, G, H on every worksheet:
#   new A (Description) = old G
#   new G (geneID)       = old H
#   new H (Count)        = old A
# Columns B..F (GeneRatio, pvalue, p.adjust, logFDR, qvalue) stay put.
# Row 1 (header) follows the same rotation: Count/...Description/geneID
#   -> Description/...geneID/Count

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 1; $r -le $lastRow; $r++) {
        $cellA = $ws.Cells.Item($r, 1)
        $cellG = $ws.Cells.Item($r, 7)
        $cellH = $ws.Cells.Item($r, 8)

        $oldA = $cellA.Value()
        $oldG = $cellG.Value()
        $oldH = $cellH.Value()

        $cellA.Value = $oldG
        $cellG.Value = $oldH
        $cellH.Value = $oldA

        # Re-assigning multi-line Description text into a different column
        # makes this host auto-expand the row height (ht=/customHeight=1).
        # The source rows never had an explicit height, so put it back to
        # the sheet's standard height and drop the custom-height flag.
        $ws.Rows.Item($r).AutoFit()
    }
}
